$wb = $excel.ActiveWorkbook

$wsRegister = $wb.Worksheets.Item("Register")
$wsLogIn = $wb.Worksheets.Item("LogIn")

# Update the shared account-name string used by both the Register sheet
# (I2) and the LogIn sheet (A2) so the two cells stay in sync and the
# old value is fully replaced.
$wsRegister.Range("I2").Value = "newAccount322788"
$wsLogIn.Range("A2").Value = "newAccount322788"

# Move the selection on the LogIn sheet (no longer the active tab).
$wsLogIn.Activate()
$wsLogIn.Range("A8").Select()

# Make Register the active tab and move its selection too.
$wsRegister.Activate()
$wsRegister.Range("G10").Select()
